$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171, shifting existing rows 171-195 down to 172-196
$ws.Range("A171").EntireRow.Insert()

# Populate the new row 171 with the inserted record
$ws.Range("A171").Value = 5
$ws.Range("B171").Value = "Macroferia Regional de Talca"
$ws.Range("C171").Value = "Maule"
$ws.Range("D171").Value = 44984
$ws.Range("E171").Value = 7
$ws.Range("F171").Value = 100112030
$ws.Range("G171").Value = "Poroto granado"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 400
$ws.Range("K171").Value = 24000
$ws.Range("L171").Value = 24000
$ws.Range("M171").Value = 24000
$ws.Range("N171").Value = "$/saco 25 kilos"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 960
$ws.Range("Q171").Value = 25
$ws.Range("R171").Value = "Hortaliza"
